# Update handback report timestamps for the 436f300c... entry (3rd data row)
# on both the zh-cn and de-de language sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn sheet: row 4 -> Correspond Handoff Datetime (D4) & Correspond Handback DateTime (G4)
$wsZhCn.Range("D4").Value = "2016-02-19 07:54:50"
$wsZhCn.Range("G4").Value = "2016-02-19 07:55:52"

# de-de sheet: row 4 -> Correspond Handoff Datetime (D4) & Correspond Handback DateTime (G4)
$wsDeDe.Range("D4").Value = "2016-02-19 07:55:05"
$wsDeDe.Range("G4").Value = "2016-02-19 07:56:16"
